# The deck ships two DrawingML theme parts:
#   theme1.xml - backs the slide master (rId12) and currently holds the
#                "Integral" design's color scheme.
#   theme2.xml - backs the notes master and currently holds the stock
#                "Office Theme" color scheme.
#
# The authored change swaps the two palettes: the presentation's live theme
# becomes the default "Office Theme" colors (what used to live in
# theme2.xml), while the old "Integral" colors end up parked in theme2.xml.
#
# PowerPoint's automation surface only ever hands back a single shared
# Theme/ThemeColorScheme object (SlideMaster.Theme, NotesMaster.Theme,
# HandoutMaster.Theme and Slide.ThemeColorScheme all resolve to the same
# live theme), so we recolor that shared ThemeColorScheme to match the
# Office Theme palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function HexToBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Item order: ppDarkColor1, ppLightColor1, ppDarkColor2, ppLightColor2,
# ppAccentColor1..6, ppHyperlinkColor, ppFollowedHyperlinkColor
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToBgr($officeThemeColors[$i - 1])
}
